# edit.ps1
# Applies "Updated capital structure database" changes:
#  - Inserts a new row for Amarin Corporation plc (NasdaqGM:AMRN) as row 3,
#    shifting Prothena Corporation plc to row 4 and Alkermes plc to row 5.
#  - Refreshes all computed capital-structure metrics for the Ireland summary
#    row and for each company row to reflect the updated dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 for the new company (Amarin Corporation plc).
# This pushes the former row 3 (Prothena) to row 4, and former row 4 (Alkermes) to row 5.
$ws.Rows.Item(3).Insert()

# --- Row 2 ---
# B2 holds a numeric-looking label ("3") that must remain stored as text
# (matching its original inline-string type), so force text formatting,
# assign the value, then clear the formatting override we just added.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").ClearFormats()
$ws.Range("D2").Value = 0.128
$ws.Range("F2").Value = 0.24
$ws.Range("G2").Value = 0.2283316187656361
$ws.Range("H2").Value = -0.1074134578338343
$ws.Range("I2").Value = -0.1136060106594885
$ws.Range("J2").Value = -0.1136060106594885
$ws.Range("K2").Value = -191.5
$ws.Range("L2").Value = -0.1086617917336464
$ws.Range("M2").Value = 37.89
$ws.Range("N2").Value = 0.006819777173815224
$ws.Range("O2").Value = -0.1978590078328982
$ws.Range("S2").Value = 37.89
$ws.Range("U2").Value = 762.6
$ws.Range("V2").Value = 0.1372594899116255
$ws.Range("W2").Value = -0.06890095487736378
$ws.Range("X2").Value = 0.06522993533495799
$ws.Range("Y2").Value = -0.1341308902123218
$ws.Range("Z2").Value = 1.904280353596428
$ws.Range("AA2").Value = 1.337005804844218
$ws.Range("AB2").Value = 0.06357235218884756
$ws.Range("AC2").Value = 1.273433452655371
$ws.Range("AD2").Value = 284.97
$ws.Range("AE2").Value = 52.26719639869408
$ws.Range("AF2").Value = 337.2371963986941
$ws.Range("AG2").Value = -425.3628036013059
$ws.Range("AH2").Value = 0.05722541070395923
$ws.Range("AI2").Value = 0.1506574305239645
$ws.Range("AJ2").Value = -0.08290804399583794
$ws.Range("AK2").Value = -0.2882179719004691
$ws.Range("AL2").Value = 13.87
$ws.Range("AM2").Value = -5.599999999999998
$ws.Range("AN2").Value = -2.730905606133206
$ws.Range("AO2").Value = -14.79452054794521
$ws.Range("AP2").Value = 4.076308611416444
$ws.Range("AQ2").Value = 36.64285714285716

# --- Row 3 ---
$ws.Range("A3").Value = "Ireland"
$ws.Range("B3").Value = "Amarin Corporation plc (NasdaqGM:AMRN)"
$ws.Range("C3").Value = "Drugs (Biotechnology)"
$ws.Range("D3").Value = 0.525
$ws.Range("G3").Value = 0.0306727673275716
$ws.Range("H3").Value = -0.03965429588205389
$ws.Range("I3").Value = -0.04054281901790672
$ws.Range("J3").Value = -0.04054281901790672
$ws.Range("K3").Value = -15.9
$ws.Range("L3").Value = -0.02694458566344687
$ws.Range("M3").Value = 30.2
$ws.Range("N3").Value = 0.0158830335542232
$ws.Range("O3").Value = -1.89937106918239
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 30.2
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 207.2
$ws.Range("V3").Value = 0.1089723361733459
$ws.Range("W3").Value = -0.02770034843205575
$ws.Range("X3").Value = 0.06357584909499579
$ws.Range("Y3").Value = -0.09127619752705154
$ws.Range("Z3").Value = -264.8073475023409
$ws.Range("AA3").Value = 10.73603636439934
$ws.Range("AB3").Value = 0.06313694853745443
$ws.Range("AC3").Value = 10.67289941586189
$ws.Range("AD3").Value = 9.470000000000001
$ws.Range("AE3").Value = 15.07158751233379
$ws.Range("AF3").Value = 24.54158751233379
$ws.Range("AG3").Value = -182.6584124876662
$ws.Range("AH3").Value = 0.0127426437392804
$ws.Range("AI3").Value = 0.03884136150163546
$ws.Range("AJ3").Value = -0.1062745056120052
$ws.Range("AK3").Value = -0.4301472532582807
$ws.Range("AL3").Value = 3.88
$ws.Range("AM3").Value = -3.54
$ws.Range("AN3").Value = -0.4955520669806385
$ws.Range("AO3").Value = -5.489690721649485
$ws.Range("AP3").Value = 9.558263343153648
$ws.Range("AQ3").Value = 6.016949152542373

# --- Row 4 ---
$ws.Range("D4").Value = -0.257
$ws.Range("G4").Value = -19.49265687583446
$ws.Range("H4").Value = -112.4165554072096
$ws.Range("I4").Value = -139.2732105882313
$ws.Range("J4").Value = -139.2732105882313
$ws.Range("K4").Value = -102
$ws.Range("L4").Value = -136.1815754339119
$ws.Range("U4").Value = 314.5
$ws.Range("V4").Value = 0.6558915537017727
$ws.Range("W4").Value = -0.3526970954356847
$ws.Range("X4").Value = 0.06522993533495799
$ws.Range("Y4").Value = -0.4179270307706426
$ws.Range("Z4").Value = -0.009599877817114071
$ws.Range("AA4").Value = 1.337005804844218
$ws.Range("AB4").Value = 0.06357235218884756
$ws.Range("AC4").Value = 1.273433452655371
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 23.07817365292621
$ws.Range("AF4").Value = 23.07817365292621
$ws.Range("AG4").Value = -291.4218263470738
$ws.Range("AH4").Value = 0.04591957005451591
$ws.Range("AI4").Value = 0.09922759857666398
$ws.Range("AJ4").Value = -1.54947180040601
$ws.Range("AK4").Value = 3.557316033366037
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = -2.93
$ws.Range("AN4").Value = -0
$ws.Range("AP4").Value = 2.970660819032353
$ws.Range("AQ4").Value = 36.07508532423208
$ws.Range("AO4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").Value = 0.128
$ws.Range("F5").Value = 0.24
$ws.Range("G5").Value = 0.3405036278275715
$ws.Range("H5").Value = -0.06973965002134017
$ws.Range("I5").Value = -0.06143703546452141
$ws.Range("J5").Value = -0.06143703546452141
$ws.Range("K5").Value = -73.59999999999999
$ws.Range("L5").Value = -0.06282543747332479
$ws.Range("M5").Value = 7.69
$ws.Range("N5").Value = 0.002422047244094488
$ws.Range("O5").Value = -0.1044836956521739
$ws.Range("S5").Value = 7.69
$ws.Range("U5").Value = 240.9
$ws.Range("V5").Value = 0.0758740157480315
$ws.Range("W5").Value = -0.06890095487736378
$ws.Range("X5").Value = 0.06725341106558885
$ws.Range("Y5").Value = -0.1361543659429526
$ws.Range("Z5").Value = 1.164840102158601
$ws.Range("AA5").Value = -0.07156432266681473
$ws.Range("AB5").Value = 0.06534949529232309
$ws.Range("AC5").Value = -0.1369138179591378
$ws.Range("AD5").Value = 275.5
$ws.Range("AE5").Value = 14.11743523343408
$ws.Range("AF5").Value = 289.6174352334341
$ws.Range("AG5").Value = 48.71743523343409
$ws.Range("AH5").Value = 0.08359290474271965
$ws.Range("AI5").Value = 0.2107814848682982
$ws.Range("AJ5").Value = 0.01511219150319432
$ws.Range("AK5").Value = 0.04299416258068413
$ws.Range("AL5").Value = 9.99
$ws.Range("AM5").Value = 0.870000000000001
$ws.Range("AN5").Value = 21.42301710730948
$ws.Range("AO5").Value = -7.827827827827828
$ws.Range("AP5").Value = 3.788292008820691
$ws.Range("AQ5").Value = -89.88505747126428
